# Auto-generated script applying the scheduled-runner data refresh
# to the Leve profit tables across all craft sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 15
$ws.Range("H15").Value = 531.3171
$ws.Range("I15").Value = 531.3171
$ws.Range("K15").Value = 1593.9513
$ws.Range("M15").Value = -1424.9513
# Row 51
$ws.Range("H51").Value = 2649.5
$ws.Range("J51").Value = 3000
$ws.Range("L51").Value = 3000
$ws.Range("N51").Value = -3968
# Row 64
$ws.Range("H64").Value = 8780.684999999999
$ws.Range("J64").Value = 9766.666999999999
$ws.Range("L64").Value = 9766.666999999999
$ws.Range("N64").Value = -10262.667
# Row 67
$ws.Range("H67").Value = 8780.684999999999
$ws.Range("J67").Value = 9766.666999999999
$ws.Range("L67").Value = 9766.666999999999
$ws.Range("N67").Value = -11482.667
# Row 132
$ws.Range("H132").Value = 1244.8
$ws.Range("I132").Value = 932.1111
$ws.Range("J132").Value = 4059
$ws.Range("K132").Value = 2796.3333
$ws.Range("L132").Value = 12177
$ws.Range("M132").Value = -266.3332999999998
$ws.Range("N132").Value = -17237
# Row 135
$ws.Range("H135").Value = 777.3421
$ws.Range("I135").Value = 815.08826
$ws.Range("J135").Value = 456.5
$ws.Range("K135").Value = 7335.79434
$ws.Range("L135").Value = 4108.5
$ws.Range("M135").Value = -4800.79434
$ws.Range("N135").Value = -9178.5
# Row 137
$ws.Range("H137").Value = 1768.625
$ws.Range("I137").Value = 1601.0769
$ws.Range("J137").Value = 2494.6667
$ws.Range("K137").Value = 4803.2307
$ws.Range("L137").Value = 7484.000100000001
$ws.Range("M137").Value = -2253.2307
$ws.Range("N137").Value = -12584.0001
# Row 138
$ws.Range("H138").Value = 3031.1191
$ws.Range("I138").Value = 1131.826
$ws.Range("J138").Value = 5330.263
$ws.Range("K138").Value = 3395.478
$ws.Range("L138").Value = 15990.789
$ws.Range("M138").Value = 1744.522
$ws.Range("N138").Value = -26270.789
# Row 141
$ws.Range("H141").Value = 2590.348
$ws.Range("I141").Value = 2572.318
$ws.Range("J141").Value = 2987
$ws.Range("K141").Value = 7716.954000000001
$ws.Range("L141").Value = 8961
$ws.Range("M141").Value = -2536.954000000001
$ws.Range("N141").Value = -19321

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 5257.6772
$ws.Range("I32").Value = 4424.5557
$ws.Range("K32").Value = 4424.5557
$ws.Range("M32").Value = -4137.5557
# Row 61
$ws.Range("H61").Value = 1787.1333
$ws.Range("I61").Value = 1539.2693
$ws.Range("K61").Value = 1539.2693
$ws.Range("M61").Value = -1327.2693
# Row 74
$ws.Range("H74").Value = 2461.8845
$ws.Range("I74").Value = 2435.1738
$ws.Range("J74").Value = 2666.6667
$ws.Range("K74").Value = 2435.1738
$ws.Range("L74").Value = 2666.6667
$ws.Range("M74").Value = -1561.1738
$ws.Range("N74").Value = -4414.6667
# Row 77
$ws.Range("H77").Value = 2461.8845
$ws.Range("I77").Value = 2435.1738
$ws.Range("J77").Value = 2666.6667
$ws.Range("K77").Value = 12175.869
$ws.Range("L77").Value = 13333.3335
$ws.Range("M77").Value = -7807.869000000001
$ws.Range("N77").Value = -22069.3335
# Row 121
$ws.Range("H121").Value = 119995
$ws.Range("J121").Value = 119995
$ws.Range("L121").Value = 119995
$ws.Range("N121").Value = -123489
# Row 122
$ws.Range("H122").Value = 2370.1516
$ws.Range("I122").Value = 2249.76
$ws.Range("K122").Value = 6749.280000000001
$ws.Range("M122").Value = -4299.280000000001
# Row 132
$ws.Range("H132").Value = 1793.0151
$ws.Range("I132").Value = 1777.705
$ws.Range("K132").Value = 5333.115
$ws.Range("M132").Value = -2803.115
# Row 135
$ws.Range("H135").Value = 263118.3
$ws.Range("J135").Value = 263118.3
$ws.Range("L135").Value = 263118.3
$ws.Range("N135").Value = -273258.3
# Row 136
$ws.Range("H136").Value = 1787.1333
$ws.Range("I136").Value = 1539.2693
$ws.Range("K136").Value = 4617.8079
$ws.Range("M136").Value = -2067.8079

$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 45458536
$ws.Range("I86").Value = 2382.4
$ws.Range("K86").Value = 2382.4
$ws.Range("M86").Value = -1259.4
# Row 89
$ws.Range("H89").Value = 45458536
$ws.Range("I89").Value = 2382.4
$ws.Range("K89").Value = 11912
$ws.Range("M89").Value = -6296
# Row 132
$ws.Range("H132").Value = 300000
$ws.Range("J132").Value = 300000
$ws.Range("L132").Value = 300000
$ws.Range("N132").Value = -310120
# Row 134
$ws.Range("H134").Value = 1536.1852
$ws.Range("I134").Value = 1207.409
$ws.Range("K134").Value = 3622.227
$ws.Range("M134").Value = -1087.227

$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 1189.0714
$ws.Range("I58").Value = 1189.0714
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 1189.0714
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -986.0714
$ws.Range("N58").Value = ""
# Row 62
$ws.Range("H62").Value = 5339.625
$ws.Range("I62").Value = 4560
$ws.Range("J62").Value = 6342
$ws.Range("K62").Value = 4560
$ws.Range("L62").Value = 6342
$ws.Range("M62").Value = -3936
$ws.Range("N62").Value = -7590
# Row 65
$ws.Range("H65").Value = 5339.625
$ws.Range("I65").Value = 4560
$ws.Range("J65").Value = 6342
$ws.Range("K65").Value = 22800
$ws.Range("L65").Value = 31710
$ws.Range("M65").Value = -19680
$ws.Range("N65").Value = -37950
# Row 132
$ws.Range("H132").Value = 1487.5869
$ws.Range("I132").Value = 1494.9697
$ws.Range("K132").Value = 4484.909100000001
$ws.Range("M132").Value = -1954.909100000001
# Row 136
$ws.Range("H136").Value = 1189.0714
$ws.Range("I136").Value = 1189.0714
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 3567.2142
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -1017.2142
$ws.Range("N136").Value = ""
# Row 137
$ws.Range("H137").Value = 186659.89
$ws.Range("J137").Value = 147498.62
$ws.Range("L137").Value = 147498.62
$ws.Range("N137").Value = -157698.62

$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Range("H2").Value = 88.75
$ws.Range("I2").Value = 32.705883
$ws.Range("J2").Value = 152.26666
$ws.Range("K2").Value = 196.235298
$ws.Range("L2").Value = 913.59996
$ws.Range("M2").Value = -83.235298
$ws.Range("N2").Value = -1139.59996
# Row 63
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").Value = ""
# Row 66
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").Value = ""
# Row 86
$ws.Range("H86").Value = 1000
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 1000
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 3000
$ws.Range("M86").Value = ""
$ws.Range("N86").Value = -5372
# Row 89
$ws.Range("H89").Value = 1000
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 1000
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 9000
$ws.Range("M89").Value = ""
$ws.Range("N89").Value = -20856
# Row 129
$ws.Range("H129").Value = 2416.3157
$ws.Range("I129").Value = 282.41666
$ws.Range("J129").Value = 6074.4287
$ws.Range("K129").Value = 847.2499799999999
$ws.Range("L129").Value = 18223.2861
$ws.Range("M129").Value = 4152.75002
$ws.Range("N129").Value = -28223.2861

$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 49.795456
$ws.Range("I2").Value = 45.19355
$ws.Range("K2").Value = 45.19355
$ws.Range("M2").Value = 67.80645
# Row 29
$ws.Range("H29").Value = 801697.6
$ws.Range("J29").Value = 1500
$ws.Range("L29").Value = 1500
$ws.Range("N29").Value = -2080
# Row 113
$ws.Range("H113").Value = 35722228
$ws.Range("I113").Value = 76928650
$ws.Range("K113").Value = 76928650
$ws.Range("M113").Value = -76926480
# Row 132
$ws.Range("H132").Value = 2253.853
$ws.Range("J132").Value = 1500
$ws.Range("L132").Value = 4500
$ws.Range("N132").Value = -9560

$ws = $wb.Worksheets.Item("LTW")
# Row 55
$ws.Range("H55").Value = 1779.6471
$ws.Range("I55").Value = 425.8
$ws.Range("J55").Value = 3713.7144
$ws.Range("K55").Value = 425.8
$ws.Range("L55").Value = 3713.7144
$ws.Range("M55").Value = -252.8
$ws.Range("N55").Value = -4059.7144
# Row 82
$ws.Range("H82").Value = 3585.7856
$ws.Range("I82").Value = 976
$ws.Range("K82").Value = 976
$ws.Range("M82").Value = -615
# Row 85
$ws.Range("H85").Value = 3585.7856
$ws.Range("I85").Value = 976
$ws.Range("K85").Value = 976
$ws.Range("M85").Value = 272
# Row 100
$ws.Range("H100").Value = 4718.1816
$ws.Range("I100").Value = 2925
$ws.Range("K100").Value = 2925
$ws.Range("M100").Value = -2384
# Row 122
$ws.Range("H122").Value = 5664.75
$ws.Range("I122").Value = 5664.75
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 16994.25
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -14544.25
$ws.Range("N122").Value = ""
# Row 132
$ws.Range("H132").Value = 2911.0408
$ws.Range("I132").Value = 1722.9286
$ws.Range("K132").Value = 5168.7858
$ws.Range("M132").Value = -2638.7858

$ws = $wb.Worksheets.Item("WVR")
# Row 64
$ws.Range("H64").Value = 95076
$ws.Range("J64").Value = 95076
$ws.Range("L64").Value = 95076
$ws.Range("N64").Value = -95572
# Row 67
$ws.Range("H67").Value = 95076
$ws.Range("J67").Value = 95076
$ws.Range("L67").Value = 95076
$ws.Range("N67").Value = -96792
# Row 122
$ws.Range("H122").Value = 4759.467
$ws.Range("I122").Value = 1941.8572
$ws.Range("J122").Value = 7224.875
$ws.Range("K122").Value = 5825.571599999999
$ws.Range("L122").Value = 21674.625
$ws.Range("M122").Value = -3375.571599999999
$ws.Range("N122").Value = -26574.625
# Row 132
$ws.Range("H132").Value = 2666.718
$ws.Range("I132").Value = 2048.5806
$ws.Range("J132").Value = 5062
$ws.Range("K132").Value = 6145.7418
$ws.Range("L132").Value = 15186
$ws.Range("M132").Value = -3615.7418
$ws.Range("N132").Value = -20246

Write-Host "Applied scheduled market-data refresh to leve profit tables."